# Adding multiple columns primary key support:
# Insert a new "data_type" column (C) that holds data types (e.g. DATETIME)
# separately from the "type" column (now D, holds PK/FK), and add a new
# "AA" table demonstrating a multi-column (composite) primary key.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts old C/D/E -> D/E/F
# and keeps their widths/values intact.
$ws.Columns("C:C").Insert()

# Add a new table "AA" with a composite (multi-column) primary key:
# columns "id" and "a" both marked PK.
$ws.Range("A13").Value = "AA"
$ws.Range("B13").Value = "id"
$ws.Range("D13").Value = "PK"

$ws.Range("A14").Value = "AA"

# New column header
$ws.Range("C1").Value = "data_type"

$ws.Range("B14").Value = "a"
$ws.Range("D14").Value = "PK"

# Move the DATETIME "type" values back into the new data_type column
# (they describe the data type of a column, not PK/FK).
$ws.Range("C3").Value = "DATETIME"
$ws.Range("D3").Value = ""

$ws.Range("C6").Value = "DATETIME"
$ws.Range("D6").Value = ""

$ws.Range("C8").Value = "DATETIME"
$ws.Range("D8").Value = ""

$ws.Range("C11").Value = "DATETIME"
$ws.Range("D11").Value = ""

# Selection / active cell as reflected in the saved view
$ws.Range("E7").Select()
